$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit reruns the NATMI pipeline with updated TPM input, which changes the
# per-pair expression/specificity statistics and narrows the result set from 9 to 8
# sending/target-cluster pairs (row 10 disappears).
$ws.Rows.Item(10).Delete()

# New row values, columns A-T (A-D are cluster/gene labels, E-T are recomputed metrics)
$newData = @(
    @("ECs", "Pspn", "Ret", "ECs", 1, 0.3333333333333333, 0.2111796666666667, 0.633539, 0.1733132136419605, 0.1733132136419605, 3, 1, 4.328495333333334, 12.985486, 0.2325244328639614, 0.2325244328639614, 0.9140902016615555, 8.226811814953999, 0.04029955670992744, 0.04029955670992743),
    @("ECs", "Pspn", "Ret", "FAPs", 1, 0.3333333333333333, 0.2111796666666667, 0.633539, 0.1733132136419605, 0.1733132136419605, 3, 1, 7.361539333333333, 22.084618, 0.3954579193622197, 0.3954579193622196, 1.554607422566889, 13.991466803102, 0.06853808286482957, 0.06853808286482956),
    @("ECs", "Pspn", "Ret", "MuSCs", 1, 0.3333333333333333, 0.2111796666666667, 0.633539, 0.1733132136419605, 0.1733132136419605, 3, 1, 6.912382333333333, 20.737147, 0.3713294477689628, 0.3713294477689628, 1.459754597025889, 13.137791373233, 0.06435629991273346, 0.06435629991273345),
    @("ECs", "Pspn", "Ret", "Resolving-Mac", 1, 0.3333333333333333, 0.2111796666666667, 0.633539, 0.1733132136419605, 0.1733132136419605, 1, 0.3333333333333333, 0.012811, 0.038433, 0.0006882000048562394, 0.0006882000048562392, 0.002705422709666667, 0.024348804387, 0.0001192741544700477, 0.0001192741544700476),
    @("FAPs", "Pspn", "Ret", "ECs", 2, 0.6666666666666666, 1.007306, 3.021918, 0.8266867863580396, 0.8266867863580396, 3, 1, 4.328495333333334, 12.985486, 0.2325244328639614, 0.2325244328639614, 4.360119320238667, 39.241073882148, 0.192224876154034, 0.1922248761540339),
    @("FAPs", "Pspn", "Ret", "FAPs", 2, 0.6666666666666666, 1.007306, 3.021918, 0.8266867863580396, 0.8266867863580396, 3, 1, 7.361539333333333, 22.084618, 0.3954579193622197, 0.3954579193622196, 7.415322739702667, 66.73790465732401, 0.3269198364973901, 0.3269198364973901),
    @("FAPs", "Pspn", "Ret", "MuSCs", 2, 0.6666666666666666, 1.007306, 3.021918, 0.8266867863580396, 0.8266867863580396, 3, 1, 6.912382333333333, 20.737147, 0.3713294477689628, 0.3713294477689628, 6.962884198660667, 62.66595778794601, 0.3069731478562294, 0.3069731478562294),
    @("FAPs", "Pspn", "Ret", "Resolving-Mac", 2, 0.6666666666666666, 1.007306, 3.021918, 0.8266867863580396, 0.8266867863580396, 1, 0.3333333333333333, 0.012811, 0.038433, 0.0006882000048562394, 0.0006882000048562392, 0.012904597166, 0.116141374494, 0.0005689258503861918, 0.0005689258503861917),
)

$colLetters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $newData[$i]
    for ($c = 0; $c -lt $colLetters.Count; $c++) {
        $ws.Range($colLetters[$c] + $rowNum).Value = $rowValues[$c]
    }
}